# Append a new scraped listing (2025-12-09 12:39 JST) to the "ランサーズ" sheet.
# - refresh the 取得日時 (scrape timestamp) on every existing data row
# - insert one new row at row 12 for the new listing, pushing the two
#   trailing rows (old 12, old 13) down to rows 13 and 14
# - rebuild the F-column hyperlinks so they line up with the shifted rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "ランサーズ" — the only sheet with per-listing rows

$newTimestamp = "2025-12-09 12:39:44"

# --- 1. Insert a new row before the current row 12 --------------------------
# This shifts old rows 12..13 down to 13..14 and grows the used range to H14.
$ws.Rows.Item(12).Insert()

# --- 2. Drop the (now stale) hyperlink bookkeeping; it will be rebuilt below.
$ws.Hyperlinks.Delete()

# --- 3. Fill in the freshly inserted row 12 with the new listing -----------
$ws.Cells.Item(12, 1).Value = $newTimestamp
$ws.Cells.Item(12, 2).Value = "自社カレンダーとGoogleカレンダーの連携エキスパート募集"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5450296"
$ws.Cells.Item(12, 7).Value = 13

# --- 4. Refresh the 取得日時 timestamp on every other data row (2..14) ------
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 5. Rebuild the F2:F14 hyperlinks in row order --------------------------
$urls = @{
    2  = "https://www.lancers.jp/work/detail/5450158"
    3  = "https://www.lancers.jp/work/detail/5450024"
    4  = "https://www.lancers.jp/work/detail/5217096"
    5  = "https://www.lancers.jp/work/detail/5442448"
    6  = "https://www.lancers.jp/work/detail/5449973"
    7  = "https://www.lancers.jp/work/detail/5449939"
    8  = "https://www.lancers.jp/work/detail/5450139"
    9  = "https://www.lancers.jp/work/detail/5449999"
    10 = "https://www.lancers.jp/work/detail/5440861"
    11 = "https://www.lancers.jp/work/detail/5449609"
    12 = "https://www.lancers.jp/work/detail/5450296"
    13 = "https://www.lancers.jp/work/detail/5449817"
    14 = "https://www.lancers.jp/work/detail/5449948"
}

for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    # Hyperlinks.Add mints a brand-new cell style; put the original
    # "Hyperlink" cell style back so formatting matches the rest of the column.
    $cell.Style = "Hyperlink"
}
